# "tela de edição funcional"
# Rows 35-40 of Sheet1 hold a little test/demo block of data. This edit
# effectively removes the old row 35 (its values ripple up into 35-39)
# and appends a brand-new row of data at the bottom (row 40).
#
# All of these cells store their contents as TEXT even though many values
# look numeric ("1", "2", "89", ...), so every write below is forced to
# stay text by using Excel's leading-apostrophe (quote-prefix) convention,
# then the style is put back to "Normal" so no stray number-format/style
# is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$text)
    $ws.Range($cell).Value = "'" + $text
    $ws.Range($cell).Style = "Normal"
}

# Row 35 <- old row 36 (date row now starts here)
Set-TextValue "A35" "01/03/2025"
Set-TextValue "B35" "1"
Set-TextValue "C35" "2"
Set-TextValue "D35" "3"
Set-TextValue "E35" "4"
Set-TextValue "F35" "1"
Set-TextValue "G35" "5"
Set-TextValue "H35" "6"
Set-TextValue "I35" "7"
Set-TextValue "J35" "89"

# Row 36 <- old row 37
Set-TextValue "A36" "01/03/2025"
Set-TextValue "B36" "1"
Set-TextValue "C36" "1"
Set-TextValue "D36" "1"
Set-TextValue "E36" "1"
Set-TextValue "F36" "0"
Set-TextValue "G36" "1"
Set-TextValue "H36" "1"
Set-TextValue "I36" "1"
Set-TextValue "J36" "admin"

# Row 37 <- old row 38
Set-TextValue "A37" "01/03/2025"
Set-TextValue "B37" "2"
Set-TextValue "C37" "2"
Set-TextValue "D37" "2"
Set-TextValue "E37" "2"
Set-TextValue "F37" "0"
Set-TextValue "G37" "2"
Set-TextValue "H37" "2"
Set-TextValue "I37" "2"
Set-TextValue "J37" "admin"

# Row 38 <- old row 39
Set-TextValue "A38" "01/03/2025"
Set-TextValue "B38" "1"
Set-TextValue "C38" "1"
Set-TextValue "D38" "1"
Set-TextValue "E38" "1"
Set-TextValue "F38" "0"
Set-TextValue "G38" "1"
Set-TextValue "H38" "1"
Set-TextValue "I38" "1"
Set-TextValue "J38" "user"

# Row 39 <- old row 40
Set-TextValue "A39" "01/03/2025"
Set-TextValue "B39" "2"
Set-TextValue "C39" "2"
Set-TextValue "D39" "2"
Set-TextValue "E39" "2"
Set-TextValue "F39" "0"
Set-TextValue "G39" "2"
Set-TextValue "H39" "2"
Set-TextValue "I39" "2"
Set-TextValue "J39" "user"

# Row 40 <- brand-new appended row
Set-TextValue "A40" "1"
Set-TextValue "B40" "2"
Set-TextValue "C40" "3"
Set-TextValue "D40" "8749"
Set-TextValue "E40" "49684"
Set-TextValue "F40" "40935"
Set-TextValue "G40" "6"
Set-TextValue "H40" "7"
Set-TextValue "I40" "8"
Set-TextValue "J40" "9"
